# Slide 25 ("Suggested task (Tuesday, Nov. 16.)"), shape 2
# ("Content Placeholder 2"): give the placeholder an explicit
# position/size, and tweak the wording of the task description.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$sh = $s.Shapes.Item(2)

# The placeholder previously had no explicit <a:xfrm> (<p:spPr/>); set an
# explicit position/size matching the target layout.
$sh.Left = 37.6669
$sh.Top = 130.16662
$sh.Width = 866.8769
$sh.Height = 374.16671

$tr = $sh.TextFrame.TextRange

# First paragraph: reword the question. Replace through a Characters()
# sub-range (rather than re-assigning Paragraphs(n).Text) so PowerPoint
# keeps a single run instead of diff-splitting it into several runs.
$para1 = $tr.Paragraphs(1)
$c1 = $tr.Characters($para1.Start, $para1.Length)
$c1.Text = "How would you express and address (mitigate) your project’s particular ethical dilemmas according to each theory? (one or two sentences for each)"

# Sixth paragraph (sub-bullet under "Care Ethics/Social Ethics of
# Engineering"): reword.
$para6 = $tr.Paragraphs(6)
$c6 = $tr.Characters($para6.Start, $para6.Length)
$c6.Text = "For the Social Ethics of Engineering, please associate your comments with the answers from the previous tasks. Add references to your claims or definitions."

# Remove the trailing empty paragraph that followed it.
$para7 = $tr.Paragraphs(7)
$para7.Delete()
